$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update temperature column (B2:B16) with new randomized "homer model" values
$values = @(47, 33, 45, 37, 34, 45, 78, 41, 68, 76, 54, 38, 75, 75, 45)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Update the selected/active cell on the sheet (active cell C10)
$ws.Range("C10").Select()
